$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.995.61"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.883.38"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.41"
$ws.Range("E5").Value = "  -3.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9990"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4924"
$ws.Range("E7").Value = "  -3.37%  "
$ws.Range("E8").Value = "  -2.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06637"
$ws.Range("E9").Value = "  -2.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.876.98"
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("E11").Value = "  -3.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07179"
$ws.Range("E12").Value = "  -1.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6687"
$ws.Range("E13").Value = "  -3.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.40"
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.898"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.944.23"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007839"
$ws.Range("E17").Value = "  -4.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9989"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.81"
$ws.Range("E19").Value = "  -2.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.118.41"
$ws.Range("E20").Value = "  -1.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9985"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.786"
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.878"
$ws.Range("E23").Value = "  +2.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.115"
$ws.Range("E24").Value = "  -1.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.84"
$ws.Range("E25").Value = "  +2.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.64"
$ws.Range("E26").Value = "  +4.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.05"
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.921"
$ws.Range("E28").Value = "  -3.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.383"
$ws.Range("E29").Value = "  -1.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.212"
$ws.Range("E30").Value = "  -1.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08760"
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.003"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05050"
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7175"
$ws.Range("E34").Value = "  -0.55%  "
$ws.Range("E35").Value = "  -1.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.668"
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01788"
$ws.Range("E37").Value = "  +5.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.698"
$ws.Range("E38").Value = "  -4.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.170"
$ws.Range("E39").Value = "  -4.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9366"
$ws.Range("E40").Value = "  -2.94%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4244"
$ws.Range("E41").Value = "  -1.62%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "103.90"
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9990"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.751"
$ws.Range("E44").Value = "  -6.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.394"
$ws.Range("E45").Value = "  -2.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1271"
$ws.Range("E46").Value = "  -0.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05706"
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "32.68"
$ws.Range("E48").Value = "  -1.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.309"
$ws.Range("E49").Value = "  -0.97%  "
$ws.Range("E50").Value = "  -1.25%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.343"
$ws.Range("E51").Value = "  -2.09%  "
